$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 187
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 53222
$ws.Range("H2").Value = 0.15

# Row 3 updates
$ws.Range("E3").Value = 110
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 1059235.75
$ws.Range("H3").Value = 9.56
